## "Generate Report for Archive"
##
## Two changes, applied via Excel COM automation against the already-open
## ActiveWorkbook:
##
## 1. Every "Status" cell that currently reads "Ready for handoff" is moved
##    on to "In Translation" (Overview!E2:F2 and the Status column - column C -
##    on the per-language report sheets zh-cn / de-de).
##
## 2. The now-narrower "Status" columns are re-sized down (they no longer need
##    to fit "Ready for handoff"): Overview columns E:F, and column C on the
##    zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- 1. Update the status text everywhere it appears -----------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($rowIdx = 1; $rowIdx -le $rowCount; $rowIdx++) {
        for ($colIdx = 1; $colIdx -le $colCount; $colIdx++) {
            $cell = $ws.Cells.Item($rowIdx, $colIdx)
            $cellText = [string]$cell.Text
            if ($cellText -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- 2. Shrink the Status columns to match the shorter text ----------------
# ColumnWidth is expressed in characters (Calibri 11 / MDW-7 grid); 12.5
# lands on the same pixel-quantized column run as the narrower report width.
$newColumnWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E:F").ColumnWidth = $newColumnWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
